$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" column (H) ---
# Copy the header formatting from G1 (bold/bordered style) onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Block 1 (Iterations = 100), rows 2-11 : refreshed D/E/F values + new H (Label) values
$ws.Range("D2").Value = 0.4977162019837889
$ws.Range("E2").Value = 0.4977162019837889
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.649636973766073
$ws.Range("E3").Value = 0.649636973766073
$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.000000413095979341948192060756
$ws.Range("E4").Value = 0.000000413095979341948192060756
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 0.5881583923622591
$ws.Range("E5").Value = 0.5881583923622591
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.5168972778802406
$ws.Range("E6").Value = 0.5168972778802406
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.6804532612688184
$ws.Range("E7").Value = 0.3195467387311816
$ws.Range("H7").Value = 1

$ws.Range("D8").Value = 0.5526585994521813
$ws.Range("E8").Value = 0.4473414005478187
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.6036585675173216
$ws.Range("E9").Value = 0.3963414324826784
$ws.Range("H9").Value = 1

$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.00000000000000055620981192364673871
$ws.Range("F11").Value = 1507.5498046875
$ws.Range("H11").Value = 1

# Block 2 (Iterations = 200), rows 12-21 : only the new H (Label) values are added
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
